# Review 167 -> Review 166 document update
# Applies the OOXML diff by editing paragraph ranges through the Word
# COM object model (ActiveDocument.Paragraphs / Range).

$d = $word.ActiveDocument

function Set-ParaText([int]$index, [string]$text) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
    $r.Text = $text
    return $p
}

# 1. Heading: "Review 167: ..." -> "Review 166: ..."
Set-ParaText 1 'Review 166: In-Context Pretraining: Language Modeling Beyond Document Boundaries' | Out-Null

# 2. Bold "Paper:" line -> new arXiv id/version
Set-ParaText 2 'Paper: https://arxiv.org/abs/2310.10638v6' | Out-Null

# 3. Huggingface link paragraph: new id + trailing line break in the same run
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4.MoveEnd(1, -1) | Out-Null
$r4.Text = "https://huggingface.co/papers/2310.10638`v"

# 4. First Hebrew intro paragraph -> new intro text
Set-ParaText 6 'כשאתם מאמנים מודל שפה (אימון מקדים) המשימה היא חיזוי הטוקן הבא. כאשר מאמנים מודל שפה בעל חלון הקשר (context) ארוך משרשרים כמה מסמכים שנבחרו באקראי ומאמנים תוך כדי חיזוי הטוקן הבא. ' | Out-Null

# 5. Paragraphs 8 and 9 ("המאמר שנסקור..." + "כלומר עבור כל טוקן...")
#    collapse into a single paragraph with new text: delete paragraph 9
#    (its Range.Delete also removes its paragraph mark, merging it into
#    paragraph 8), then overwrite paragraph 8's text.
$d.Paragraphs.Item(9).Range.Delete() | Out-Null
Set-ParaText 8 'המאמר שנסקור היום ב-#shorthebrewpapereviews משכלל את הגישה הזו ומציע לשרשר מסמכים שהם קרובים מבחינת המשמעות אחד לשני במקום לבחור אותם באקראי. איך נבחרים מסמכים קרובים – לפי המרחק בין השיכונים(embedding) שלהם. אבל יש בעיה קטנה עם הגישה הנאיבית הזו. יש מסמכים שהם דומים ליותר מדי מסמכים ואז המודל "יראה״ אותם יותר פעמים מהאחרים שעלול כמובן לפגוע בביצועי המודל המאומן (יוצר סוג של overfit). ' | Out-Null

# 6. Paragraph 10 ("מודל התגמול מאומן...") -> new "כדי להתגבר..." text
Set-ParaText 10 'כדי להתגבר על סוגיה זו המחברים מציעים לתאר את כל המסמכים בדאטהסט על ידי גרף שמשקל של כל קשת בו (בין שני המסמכים) שווה לדמיון ביניהם.ֿ אחרי שיש לנו ביד גרף כזה ניתן לתאר את הבעיה בתור בעייה דומה לזו של איש מכירות המטייל (maximum travelling salesman problem) כאשר המטרה כאן למצוא מסלולים זרים (שהאיחוד שלהם מכיל את כל הקודקודים וכל קודקוד מופיע רק פעם אחת באיחוד הזה). פותרים את הבעיה הזו עם אלגוריתם די אינטואיטיבי. ' | Out-Null

# 7. After paragraph 11 (the empty paragraph following the paragraph from
#    step 6), insert three new paragraphs: text, empty, text.
$d.Paragraphs.Item(11).Range.InsertParagraphAfter() | Out-Null
Set-ParaText 12 'לקודקוד נתון בוחרים כמה קודקודים דומים (NN-nearest neighbors) ובונים מהם מסלול בעל משקל כולל מקסימלי  (סכום של כל משקלי הקשתות). כל פעם בוחרים קודקוד (מסמך) הקרוב ביותר לקודקוד האחרון שנבחר. מספר NN בכל תת-מסלול נבחר לפי אורך הקונטקסט (אורכו של כל שרשור המסמכים שווה לאורך הקונטקסט). אחרי שמסיימים לבנות כל שרשור מורידים את קוקודיו מהגרף הכולל. ' | Out-Null

$d.Paragraphs.Item(12).Range.InsertParagraphAfter() | Out-Null
# paragraph 13 is the newly inserted empty paragraph; leave it blank

$d.Paragraphs.Item(13).Range.InsertParagraphAfter() | Out-Null
Set-ParaText 14 'לאחר מכן בוחרים מסמך עם הדרגה הכוללת המינימלית (השווה לסכום משקלי הקשתות שיוצאות ממנו) וחוזרים על התהליך. כך גורמים לכל מסמך להיכנס לשרשור עם מסמכים שכמה שיותר דומים לו.' | Out-Null
